$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Tocantins -> Rondônia
$ws.Range("A2").Value = "Rondônia"
$ws.Range("B2").Value = "Diferença 2025/07 - 2025/07"
$ws.Range("C2").Value = 1.24

# Row 3: Piauí -> Amapá
$ws.Range("A3").Value = "Amapá"
$ws.Range("B3").Value = "Diferença 2025/07 - 2025/07"
$ws.Range("C3").Value = 1.05

# Row 4: Acre (name unchanged)
$ws.Range("B4").Value = "Diferença 2025/07 - 2025/07"
$ws.Range("C4").Value = 1

# Row 5: Amazonas -> Mato Grosso
$ws.Range("A5").Value = "Mato Grosso"
$ws.Range("B5").Value = "Diferença 2025/07 - 2025/07"
$ws.Range("C5").Value = 1
$ws.Range("D5").Value = "3º"

# Row 6: Pará -> Alagoas
$ws.Range("A6").Value = "Alagoas"
$ws.Range("B6").Value = "Diferença 2025/07 - 2025/07"
$ws.Range("C6").Value = 0.99

# Row 7: Alagoas -> Ceará
$ws.Range("A7").Value = "Ceará"
$ws.Range("B7").Value = "Diferença 2025/07 - 2025/07"
$ws.Range("C7").Value = 0.96

# Row 8: Sergipe (name unchanged)
$ws.Range("B8").Value = "Diferença 2025/07 - 2025/07"
$ws.Range("C8").Value = 0.92
$ws.Range("D8").Value = "11º"

# Row 9: Brasil (name unchanged)
$ws.Range("B9").Value = "Diferença 2025/07 - 2025/07"
$ws.Range("C9").Value = 0.87

# Row 10: Nordeste (name unchanged)
$ws.Range("B10").Value = "Diferença 2025/07 - 2025/07"
$ws.Range("C10").Value = 0.9
